# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.332.28"
$ws.Range("E2").Value = "  +2.25%  "
$ws.Range("D3").Value = "'2.424.57"
$ws.Range("E3").Value = "  +3.16%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'556.99"
$ws.Range("E5").Value = "  +2.23%  "
$ws.Range("E6").Value = "  +5.00%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +2.03%  "
$ws.Range("D9").Value = "'2.425.16"
$ws.Range("E9").Value = "  +3.14%  "
$ws.Range("E10").Value = "  +5.39%  "
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("E12").Value = "  +1.65%  "
$ws.Range("E13").Value = "  +3.10%  "
$ws.Range("E14").Value = "  +6.44%  "
$ws.Range("E15").Value = "  +9.67%  "
$ws.Range("D16").Value = "'2.866.27"
$ws.Range("E16").Value = "  +3.35%  "
$ws.Range("D17").Value = "'62.118.30"
$ws.Range("E17").Value = "  +2.07%  "
$ws.Range("D18").Value = "'2.425.04"
$ws.Range("E18").Value = "  +3.43%  "
$ws.Range("E19").Value = "  +4.57%  "
$ws.Range("D20").Value = "'324.61"
$ws.Range("E20").Value = "  +1.61%  "
$ws.Range("E21").Value = "  +1.40%  "
$ws.Range("E22").Value = "  +3.25%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("E24").Value = "  +2.29%  "
$ws.Range("D25").Value = "'65.01"
$ws.Range("E25").Value = "  +2.49%  "
$ws.Range("D26").Value = "'9.20"
$ws.Range("E26").Value = "  +10.62%  "
$ws.Range("D27").Value = "'572.01"
$ws.Range("E27").Value = "  +14.21%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "'0.0₃0954"
$ws.Range("E28").Value = "  +9.80%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "'2.541.74"
$ws.Range("E29").Value = "  +3.27%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").Value = "'8.41"
$ws.Range("E31").Value = "  +5.65%  "
$ws.Range("D32").Value = "'1.46"
$ws.Range("E32").Value = "  +5.67%  "
$ws.Range("E33").Value = "  +1.55%  "
$ws.Range("E34").Value = "  +4.69%  "
$ws.Range("E35").Value = "  +4.66%  "
$ws.Range("D36").Value = "'5.75"
$ws.Range("E36").Value = "  +8.62%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  +4.63%  "
$ws.Range("D39").Value = "'0.386"
$ws.Range("E39").Value = "  +2.27%  "
$ws.Range("E40").Value = "  +4.67%  "
$ws.Range("E41").Value = "  +0.91%  "
$ws.Range("D42").Value = "'148.27"
$ws.Range("E42").Value = "  +4.23%  "
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("E44").Value = "  +2.77%  "
$ws.Range("D45").Value = "'2.33"
$ws.Range("E45").Value = "  +12.44%  "
$ws.Range("D46").Value = "'152.22"
$ws.Range("E46").Value = "  +6.85%  "
$ws.Range("E47").Value = "  +2.12%  "
$ws.Range("D48").Value = "'0.0544"
$ws.Range("E48").Value = "  +5.11%  "
$ws.Range("D49").Value = "'20.46"
$ws.Range("E49").Value = "  +7.19%  "
$ws.Range("E50").Value = "  +3.76%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0228"
$ws.Range("E51").Value = "  +3.73%  "
